$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows 239-244 (dates 2021-04-27 .. 2021-05-02), one row per
# day, columns A (date, serial) through AX (49 municipality/total counters).
$newRows = @(
    @(44313,1,1,0,15,8,0,7,1,0,0,1,2,6,0,1,0,12,2,2,6,12,0,2,0,7,1,1,0,0,1,1,9,1,4,0,0,2,8,0,2,118,0,0,0,2,0,0,0,0),
    @(44314,0,0,0,2,1,1,0,0,2,0,1,1,4,0,0,0,1,0,0,1,20,1,2,1,1,0,1,0,0,0,1,3,0,0,0,1,4,2,0,2,53,0,0,0,0,0,0,0,0),
    @(44315,5,3,0,8,5,1,2,1,1,1,6,6,8,0,0,0,2,0,3,3,34,7,5,2,4,0,0,0,0,4,3,14,6,1,0,3,4,4,0,5,153,1,0,0,0,0,0,1,0),
    @(44316,9,0,0,32,8,1,2,5,1,0,1,8,8,0,1,1,2,1,1,4,36,5,1,2,6,0,1,0,2,1,2,11,8,1,1,9,5,8,1,5,193,2,0,0,1,0,0,0,0),
    @(44317,7,0,0,16,5,4,4,2,1,1,1,4,7,0,3,0,12,2,1,7,53,3,5,2,3,0,1,0,0,2,3,19,6,2,0,4,5,9,0,7,205,0,1,0,2,0,1,0,0),
    @(44318,5,0,1,14,5,1,3,3,1,1,3,7,4,0,0,1,5,2,1,1,31,1,2,2,4,1,0,0,1,2,0,14,5,4,0,2,0,2,0,4,134,0,0,0,0,1,0,0,0)
)

$startRow = 239
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $rowVals = $newRows[$i]
    $excelRow = $startRow + $i
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($excelRow, $c + 1).Value2 = $rowVals[$c]
    }
}

# The date cells in column A use a custom date/time style (border, bold,
# centered, custom number format) that already exists on row 238 - copy
# that formatting down onto the new date cells instead of re-deriving it,
# so the same style index is reused rather than creating a duplicate.
$ws.Range("A238").Copy()
$ws.Range("A239:A244").PasteSpecial(-4122)
$excel.CutCopyMode = 0

